# "arreglado una cosa de la cookie" - reposition two placeholders that were
# slightly misaligned (detected right after the original delivery).
#
# NOTE: PowerPoint's Shape.Left/Shape.Top/Shape.Width/Shape.Height are
# expressed in points (1 pt = 12700 EMU), not EMU, so the target EMU offsets
# from the OOXML are converted to points below.

$p = $ppt.ActivePresentation

# Slide 23, shape "Marcador de contenido 2": move up slightly
# (a:off y 1517260 -> 1490627 EMU), keep its x offset and size untouched.
$slide23 = $p.Slides.Item(23)
$contentPlaceholder = $slide23.Shapes.Item(2)
$contentPlaceholder.Top = 117.37224578857422

# Slide 27, shape "Content Placeholder 22": move up and to the left a bit
# (a:off 7859485,2198914 -> 7830094,1772786 EMU), size stays the same.
$slide27 = $p.Slides.Item(27)
$qaPlaceholder = $slide27.Shapes.Item(4)
$qaPlaceholder.Left = 616.5428466796875
$qaPlaceholder.Top = 139.58949279785156
